$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "48.256.05"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.81%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.520.92"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.31"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.97"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.18%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +4.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.41"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.23"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +9.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0819"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.34%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.80%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.29"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.917.56"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.518.27"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.860"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "48.135.57"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.73%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.64"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.50%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.33%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.40"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "269.12"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +8.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.58"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.19"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.29%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.30"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.16"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.94%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.145"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +4.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.95"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.14%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.96"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.18%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0792"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.53%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.72"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.99"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.45%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.39"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +5.85%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.71%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "118.76"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.51%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.004.15"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.12%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.02%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +5.99%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.78%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.26"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.40"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.74%  "
